# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# Rebuilds the worker/period detail table (rows 16-43, columns C:F) on sheet
# "Hoja1" so that it is grouped by worker (descending by period 2311 -> 2305)
# instead of grouped by period, and updates the "Valor Mora" (F) amounts to
# match: 37333 for period 2311, 46400 for periods 2305-2310.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$data = @(
    ,@(16, "60405011", "ROSARIO MARIA FLOREZ TEHERAN", "2311", 37333)
    ,@(17, "60405011", "ROSARIO MARIA FLOREZ TEHERAN", "2310", 46400)
    ,@(18, "60405011", "ROSARIO MARIA FLOREZ TEHERAN", "2309", 46400)
    ,@(19, "60405011", "ROSARIO MARIA FLOREZ TEHERAN", "2308", 46400)
    ,@(20, "60405011", "ROSARIO MARIA FLOREZ TEHERAN", "2307", 46400)
    ,@(21, "60405011", "ROSARIO MARIA FLOREZ TEHERAN", "2306", 46400)
    ,@(22, "60405011", "ROSARIO MARIA FLOREZ TEHERAN", "2305", 46400)
    ,@(23, "1047455587", "FERNANDO ANDRES OVALLE CORDERO", "2311", 37333)
    ,@(24, "1047455587", "FERNANDO ANDRES OVALLE CORDERO", "2310", 46400)
    ,@(25, "1047455587", "FERNANDO ANDRES OVALLE CORDERO", "2309", 46400)
    ,@(26, "1047455587", "FERNANDO ANDRES OVALLE CORDERO", "2308", 46400)
    ,@(27, "1047455587", "FERNANDO ANDRES OVALLE CORDERO", "2307", 46400)
    ,@(28, "1047455587", "FERNANDO ANDRES OVALLE CORDERO", "2306", 46400)
    ,@(29, "1047455587", "FERNANDO ANDRES OVALLE CORDERO", "2305", 46400)
    ,@(30, "7920858", "RUBEN DARIO CORONEL MORALES", "2311", 37333)
    ,@(31, "7920858", "RUBEN DARIO CORONEL MORALES", "2310", 46400)
    ,@(32, "7920858", "RUBEN DARIO CORONEL MORALES", "2309", 46400)
    ,@(33, "7920858", "RUBEN DARIO CORONEL MORALES", "2308", 46400)
    ,@(34, "7920858", "RUBEN DARIO CORONEL MORALES", "2307", 46400)
    ,@(35, "7920858", "RUBEN DARIO CORONEL MORALES", "2306", 46400)
    ,@(36, "7920858", "RUBEN DARIO CORONEL MORALES", "2305", 46400)
    ,@(37, "1065003738", "JUAN MANUEL DEGIOVANNI PRECIADO", "2311", 37333)
    ,@(38, "1065003738", "JUAN MANUEL DEGIOVANNI PRECIADO", "2310", 46400)
    ,@(39, "1065003738", "JUAN MANUEL DEGIOVANNI PRECIADO", "2309", 46400)
    ,@(40, "1065003738", "JUAN MANUEL DEGIOVANNI PRECIADO", "2308", 46400)
    ,@(41, "1065003738", "JUAN MANUEL DEGIOVANNI PRECIADO", "2307", 46400)
    ,@(42, "1065003738", "JUAN MANUEL DEGIOVANNI PRECIADO", "2306", 46400)
    ,@(43, "1065003738", "JUAN MANUEL DEGIOVANNI PRECIADO", "2305", 46400)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]   # C: N° Doc Trabajador
    $ws.Cells.Item($r, 4).Value = $row[2]   # D: Nombre Trabajador
    $ws.Cells.Item($r, 5).Value = $row[3]   # E: Periodo Mora
    $ws.Cells.Item($r, 6).Value = $row[4]   # F: Valor Mora
}

$wb.Save()
